# Add analysed results for high density - reference
#
# Fills in the "Started flights" column (C) and the analysed/recalculated
# values for "High traffic density" rows 3-12 (Reference scenario), adds the
# AVERAGE / STDEV.P formulas for the new column C on the summary rows 13-14,
# and makes "High traffic density" the active sheet/tab with A15 selected
# (mirroring the previous active sheet/selection on "Medium traffic density").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("High traffic density")

# Row => Started flights(C), Completed flights(D), Avg. distance flown(E),
#        Average flying time(F), Traffic density(H), Average conflict time(J)
$data = @(
    @{ Row=3;  C=2400; D=1808; E=3735.70248969963;   F=294.68683628318598; H=10.485575676776801; J=9.5449016825461204 }
    @{ Row=4;  C=2410; D=1789; E=3822.1043566995099; F=301.06950810508698; H=10.996842018938599;  J=9.7129195574723397 }
    @{ Row=5;  C=2421; D=1817; E=3739.3780349748499; F=290.41788662630802; H=10.6583477799773;    J=10.002614116227599 }
    @{ Row=6;  C=2398; D=1818; E=3740.6256920923502; F=295.40525302530301; H=10.446655875114599;  J=9.3313289236319896 }
    @{ Row=7;  C=2409; D=1805; E=3743.6995201812301; F=291.84487534626101; H=10.5704210826934;    J=9.2148626817447497 }
    @{ Row=8;  C=2409; D=1811; E=3747.3955310451602; F=295.01653782440701; H=10.6056174365527;    J=9.0470817120622602 }
    @{ Row=9;  C=2415; D=1837; E=3788.7481489264501; F=296.023679912901;   H=10.5668559471232;    J=9.4718395461912497 }
    @{ Row=10; C=2424; D=1800; E=3839.6599836078799; F=297.17938888899999; H=10.972816612706101;  J=9.5192636173847003 }
    @{ Row=11; C=2401; D=1811; E=3811.90502405593;   F=295.50935946990597; H=10.344715637491401;  J=9.0760257441673406 }
    @{ Row=12; C=2412; D=1807; E=3745.0817024089502; F=296.35486995019397; H=10.7461598288936;    J=9.2364649681528697 }
)

foreach ($d in $data) {
    $ws.Cells.Item($d.Row, 3).Value  = $d.C   # C: Started flights (new)
    $ws.Cells.Item($d.Row, 4).Value  = $d.D   # D: Completed flights
    $ws.Cells.Item($d.Row, 5).Value  = $d.E   # E: Avg. Distance flown
    $ws.Cells.Item($d.Row, 6).Value  = $d.F   # F: Average flying time
    $ws.Cells.Item($d.Row, 8).Value  = $d.H   # H: Traffic density (new)
    $ws.Cells.Item($d.Row, 10).Value = $d.J   # J: Average conflict time
}

# Row 13 (Average) / Row 14 (Standard dev) - extend formulas to new column C
$ws.Range("C13").Formula = "=AVERAGE(C3:C12)"
$ws.Range("C13").NumberFormat = "0.00"
$ws.Range("C13").Font.Bold = $true

$ws.Range("C14").Formula = "=_xlfn.STDEV.P(C3:C12)"
$ws.Range("C14").NumberFormat = "0.00"

# "High traffic density" becomes the active/selected tab (A15 selected),
# "Medium traffic density" loses tab-selection.
$ws.Activate()
$ws.Range("A15").Select()
